$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2092.08
$ws.Range("I51").Value = 1923.0769
$ws.Range("J51").Value = 2275.1667
$ws.Range("K51").Value = 1923.0769
$ws.Range("L51").Value = 2275.1667
$ws.Range("M51").Value = -1439.0769
$ws.Range("N51").Value = -3243.1667

$ws.Range("H74").Value = 4228.5713
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 4320
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 4320
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -6192

$ws.Range("H77").Value = 4228.5713
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 4320
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 21600
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -30960

$ws.Range("H116").Value = 3519.348
$ws.Range("J116").Value = 3998.6365
$ws.Range("L116").Value = 3998.6365
$ws.Range("N116").Value = -10882.6365

$ws.Range("H127").Value = 1069.5454
$ws.Range("I127").Value = 693.5
$ws.Range("J127").Value = 2072.3333
$ws.Range("K127").Value = 2080.5
$ws.Range("L127").Value = 6216.999899999999
$ws.Range("M127").Value = 2879.5
$ws.Range("N127").Value = -16136.9999

$ws.Range("H132").Value = 3774938
$ws.Range("I132").Value = 4546390.5
$ws.Range("J132").Value = 3392
$ws.Range("K132").Value = 13639171.5
$ws.Range("L132").Value = 10176
$ws.Range("M132").Value = -13636641.5
$ws.Range("N132").Value = -15236

$ws.Range("H138").Value = 3938.2979
$ws.Range("I138").Value = 2028
$ws.Range("J138").Value = 6109.091
$ws.Range("K138").Value = 6084
$ws.Range("L138").Value = 18327.273
$ws.Range("M138").Value = -944
$ws.Range("N138").Value = -28607.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16668838
$ws.Range("I2").Value = 50001470
$ws.Range("J2").Value = 2520
$ws.Range("K2").Value = 50001470
$ws.Range("L2").Value = 2520
$ws.Range("M2").Value = -50001357
$ws.Range("N2").Value = -2746

$ws.Range("H116").Value = 16668838
$ws.Range("I116").Value = 50001470
$ws.Range("J116").Value = 2520
$ws.Range("K116").Value = 50001470
$ws.Range("L116").Value = 2520
$ws.Range("M116").Value = -49999176
$ws.Range("N116").Value = -7108

$ws.Range("H132").Value = 38466000
$ws.Range("I132").Value = 50004384
$ws.Range("J132").Value = 4733.1665
$ws.Range("K132").Value = 150013152
$ws.Range("L132").Value = 14199.4995
$ws.Range("M132").Value = -150010622
$ws.Range("N132").Value = -19259.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16668838
$ws.Range("I3").Value = 50001470
$ws.Range("J3").Value = 2520
$ws.Range("K3").Value = 50001470
$ws.Range("L3").Value = 2520
$ws.Range("M3").Value = -50001356
$ws.Range("N3").Value = -2748

$ws.Range("H22").Value = 350
$ws.Range("I22").Value = 350
$ws.Range("K22").Value = 350
$ws.Range("M22").Value = -177

$ws.Range("H53").Value = 40780
$ws.Range("J53").Value = 40780
$ws.Range("L53").Value = 40780
$ws.Range("N53").Value = -41928

$ws.Range("H134").Value = 3365.5293
$ws.Range("I134").Value = 1767.8334
$ws.Range("K134").Value = 5303.5002
$ws.Range("M134").Value = -2768.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1874.875
$ws.Range("J22").Value = 3549.75
$ws.Range("L22").Value = 3549.75
$ws.Range("N22").Value = -4249.75

$ws.Range("H99").Value = 1970.75
$ws.Range("I99").Value = 1312
$ws.Range("J99").Value = 2102.5
$ws.Range("K99").Value = 1312
$ws.Range("L99").Value = 2102.5
$ws.Range("M99").Value = 186
$ws.Range("N99").Value = -5098.5

$ws.Range("H126").Value = 1970.75
$ws.Range("I126").Value = 1312
$ws.Range("J126").Value = 2102.5
$ws.Range("K126").Value = 3936
$ws.Range("L126").Value = 6307.5
$ws.Range("M126").Value = -1466
$ws.Range("N126").Value = -11247.5

$ws.Range("H134").Value = 2852.84
$ws.Range("I134").Value = 1731.0769
$ws.Range("J134").Value = 4068.0833
$ws.Range("K134").Value = 5193.2307
$ws.Range("L134").Value = 12204.2499
$ws.Range("M134").Value = -2658.2307
$ws.Range("N134").Value = -17274.2499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1910.1666
$ws.Range("I68").Value = 799.5
$ws.Range("J68").Value = 2257.25
$ws.Range("K68").Value = 2398.5
$ws.Range("L68").Value = 6771.75
$ws.Range("M68").Value = -1587.5
$ws.Range("N68").Value = -8393.75

$ws.Range("H71").Value = 1910.1666
$ws.Range("I71").Value = 799.5
$ws.Range("J71").Value = 2257.25
$ws.Range("K71").Value = 7195.5
$ws.Range("L71").Value = 20315.25
$ws.Range("M71").Value = -3139.5
$ws.Range("N71").Value = -28427.25

$ws.Range("H107").Value = 1214.0714
$ws.Range("I107").Value = 746
$ws.Range("J107").Value = 1341.7273
$ws.Range("K107").Value = 2238
$ws.Range("L107").Value = 4025.1819
$ws.Range("M107").Value = -318
$ws.Range("N107").Value = -7865.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2769.6
$ws.Range("I113").Value = 1449.1428
$ws.Range("J113").Value = 3925
$ws.Range("K113").Value = 1449.1428
$ws.Range("L113").Value = 3925
$ws.Range("M113").Value = 720.8571999999999
$ws.Range("N113").Value = -8265

$ws.Range("H132").Value = 22730026
$ws.Range("I132").Value = 31251598
$ws.Range("K132").Value = 93754794
$ws.Range("M132").Value = -93752264

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1888.6364
$ws.Range("I22").Value = 725
$ws.Range("J22").Value = 2553.5715
$ws.Range("K22").Value = 725
$ws.Range("L22").Value = 2553.5715
$ws.Range("M22").Value = -430
$ws.Range("N22").Value = -3143.5715

$ws.Range("H27").Value = 1888.6364
$ws.Range("I27").Value = 725
$ws.Range("J27").Value = 2553.5715
$ws.Range("K27").Value = 725
$ws.Range("L27").Value = 2553.5715
$ws.Range("M27").Value = -618
$ws.Range("N27").Value = -2767.5715

$ws.Range("H68").Value = 2126.5334
$ws.Range("I68").Value = 999.6
$ws.Range("J68").Value = 2690
$ws.Range("K68").Value = 999.6
$ws.Range("L68").Value = 2690
$ws.Range("M68").Value = -250.6
$ws.Range("N68").Value = -4188

$ws.Range("H71").Value = 2126.5334
$ws.Range("I71").Value = 999.6
$ws.Range("J71").Value = 2690
$ws.Range("K71").Value = 4998
$ws.Range("L71").Value = 13450
$ws.Range("M71").Value = -1254
$ws.Range("N71").Value = -20938

$ws.Range("H132").Value = 3441.5417
$ws.Range("I132").Value = 2300.25
$ws.Range("J132").Value = 4582.8335
$ws.Range("K132").Value = 6900.75
$ws.Range("L132").Value = 13748.5005
$ws.Range("M132").Value = -4370.75
$ws.Range("N132").Value = -18808.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1999.6666
$ws.Range("I96").Value = 1999
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1999
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -626
$ws.Range("N96").Value = -4746

$ws.Range("H132").Value = 12615.615
$ws.Range("I132").Value = 4000.5715
$ws.Range("K132").Value = 12001.7145
$ws.Range("M132").Value = -9471.7145
